$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-15 (columns A=Giorno, B=TISG, C=fcs, D=buy, E=MYDIR, F=need_to_buy)
$data = @(
    @(2,  45982, 10559.6536995802, 9461.61409666624, 15411.4,  5600.33625756419,  -14.5604019070656),
    @(3,  45983, 4533.75423640896, 7279.50396665755, 9003.4,   6132.00196195993,  183.671080359062),
    @(4,  45984, 4625.19171593854, 7643.78191263307, 9003.4,   6415.21209780974,  210.649750435117),
    @(5,  45985, 12137.029246834,  10862.3031471077, 9003.4,   7217.06462204282,  378.165323714604),
    @(6,  45986, 11899.9068523764, 10508.3701784901, 9003.4,   7000.62446574863,  354.399776843279),
    @(7,  45987, 10978.1429954002, 10365.235982583,  9003.4,   6407.76374548002,  323.733322002625),
    @(8,  45988, 10978.1429954002, 10626.903869457,  9003.4,   6407.76374548002,  334.636150622375),
    @(9,  45989, 10978.1429954002, 10083.5045212309, 9003.4,   6407.76374548002,  311.994511112954),
    @(10, 45990, 4441.27964694584, 7128.65302518155, 9003.4,   6031.6631998712,   173.204842710531),
    @(11, 45991, 4285.59691763933, 6988.42367601363, 9003.4,   6023.17493479135,  167.008275450207),
    @(12, 45992, 11775.3234136769, 10796.1457573813, 8664.26,  7488.24279155423,  400.838689538981),
    @(13, 45993, 11775.3234136769, 10755.1045337058, 8664.26,  7488.24279155423,  399.128638552501),
    @(14, 45994, 11775.3234136769, 10611.1543678437, 8664.26,  7488.24279155423,  393.130714974915),
    @(15, 45995, 11775.3234136769, 10859.4356594979, 8664.26,  7488.24279155423,  403.475768793837)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
